$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new rows for May 7 - May 10, 2021, reusing the date formatting
# already used by the existing date column (copy style from A13, then
# overwrite the formula/value so no new number-format style is created).
$ws.Range("A13").Copy($ws.Range("A14:A17"))

$ws.Range("A14").Formula = "=DATE(2021,5,7)"
$ws.Range("B14").Value = "500gm"
$ws.Range("C14").Value = 45

$ws.Range("A15").Formula = "=DATE(2021,5,8)"
$ws.Range("B15").Value = "500gm"
$ws.Range("C15").Value = 45

$ws.Range("A16").Formula = "=DATE(2021,5,9)"
$ws.Range("B16").Value = "0gm"
$ws.Range("C16").Value = 0

$ws.Range("A17").Formula = "=DATE(2021,5,10)"
$ws.Range("B17").Value = "0gm"
$ws.Range("C17").Value = 0

# Update the saved selection to A2
$ws.Range("A2").Select()
